$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) == FF6495ED, matches existing "HyperLink" cell font
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) : add row 4 for 2ea98654-85a8-43f2-805f-06c50a72ef06.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.md"

$wsOverview.Range("B4").Value = "e2e\2ea98654-85a8-43f2-805f-06c50a72ef06.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5edbe2926b29fe23ceac6d752ac038d1c7867f1/e2e/2ea98654-85a8-43f2-805f-06c50a72ef06.md", [Type]::Missing, [Type]::Missing, "e2e\2ea98654-85a8-43f2-805f-06c50a72ef06.md") | Out-Null
$wsOverview.Range("B4").Font.Color = $hyperlinkColor
$wsOverview.Range("B4").Font.Underline = 2

$wsOverview.Range("C4").Value = ".md"

$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"

$wsOverview.Range("G4").Value = "2016-08-17 12:42:35"
$wsOverview.Range("G4").NumberFormat = $dateFormat

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) : add row 4
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5edbe2926b29fe23ceac6d752ac038d1c7867f1/e2e/2ea98654-85a8-43f2-805f-06c50a72ef06.md", [Type]::Missing, [Type]::Missing, "2ea98654-85a8-43f2-805f-06c50a72ef06.md") | Out-Null
$wsZhCn.Range("A4").Font.Color = $hyperlinkColor
$wsZhCn.Range("A4").Font.Underline = 2

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"

$wsZhCn.Range("G4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.012c4e42a2c7619f40e740a89eee8c8ec43f08a5.zh-cn.xlf"

$wsZhCn.Range("H4").Value = "2016-08-17 12:42:30"
$wsZhCn.Range("H4").NumberFormat = $dateFormat

$wsZhCn.Range("I4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e5edbe2926b29fe23ceac6d752ac038d1c7867f1/e2e/2ea98654-85a8-43f2-805f-06c50a72ef06.md", [Type]::Missing, [Type]::Missing, "2ea98654-85a8-43f2-805f-06c50a72ef06.md") | Out-Null
$wsZhCn.Range("I4").Font.Color = $hyperlinkColor
$wsZhCn.Range("I4").Font.Underline = 2

$wsZhCn.Range("J4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.012c4e42a2c7619f40e740a89eee8c8ec43f08a5.zh-cn.xlf"

$wsZhCn.Range("K4").Value = "2016-08-17 12:42:48"
$wsZhCn.Range("K4").NumberFormat = $dateFormat

$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) : add row 4
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5edbe2926b29fe23ceac6d752ac038d1c7867f1/e2e/2ea98654-85a8-43f2-805f-06c50a72ef06.md", [Type]::Missing, [Type]::Missing, "2ea98654-85a8-43f2-805f-06c50a72ef06.md") | Out-Null
$wsDeDe.Range("A4").Font.Color = $hyperlinkColor
$wsDeDe.Range("A4").Font.Underline = 2

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"

$wsDeDe.Range("G4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.012c4e42a2c7619f40e740a89eee8c8ec43f08a5.de-de.xlf"

$wsDeDe.Range("H4").Value = "2016-08-17 12:42:35"
$wsDeDe.Range("H4").NumberFormat = $dateFormat

$wsDeDe.Range("I4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e5edbe2926b29fe23ceac6d752ac038d1c7867f1/e2e/2ea98654-85a8-43f2-805f-06c50a72ef06.md", [Type]::Missing, [Type]::Missing, "2ea98654-85a8-43f2-805f-06c50a72ef06.md") | Out-Null
$wsDeDe.Range("I4").Font.Color = $hyperlinkColor
$wsDeDe.Range("I4").Font.Underline = 2

$wsDeDe.Range("J4").Value = "2ea98654-85a8-43f2-805f-06c50a72ef06.012c4e42a2c7619f40e740a89eee8c8ec43f08a5.de-de.xlf"

$wsDeDe.Range("K4").Value = "2016-08-17 12:42:56"
$wsDeDe.Range("K4").NumberFormat = $dateFormat

$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))
